$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing row 75 (columns B..AG)
$row75 = [ordered]@{
    "B"  = 370243
    "C"  = 50969
    "D"  = 27477
    "E"  = 23492
    "F"  = 8979
    "G"  = 854
    "H"  = 8124
    "I"  = 1142
    "J"  = 11819
    "K"  = 0
    "L"  = 11819
    "M"  = 185931
    "N"  = 164226
    "O"  = 9589
    "P"  = 12116
    "Q"  = 690
    "R"  = 690
    "S"  = 110714
    "T"  = -284141
    "U"  = 654385
    "V"  = 70514
    "W"  = 2
    "X"  = 70512
    "Y"  = 165235
    "Z"  = 5909
    "AA" = 159326
    "AB" = 305887
    "AC" = 108589
    "AD" = 197299
    "AE" = 112748
    "AF" = 6473
    "AG" = 4088
}

foreach ($col in $row75.Keys) {
    $ws.Range("$col`75").Value = $row75[$col]
}

# New row 76, starting with the series label. Column A stores these labels
# as plain text (shared strings), but assigning a date-shaped string such as
# "01-04-2021" straight to .Value gets auto-converted to a date serial by
# value-parsing heuristics. Instead, write it as a text formula (a quoted
# string literal always evaluates to genuine text) and then flatten the
# formula down to its cached value via copy / paste-special-values - this
# yields a plain text cell with no special number format, matching the rest
# of the column and leaving the style table untouched.
$ws.Range("A76").Formula = "=""01-04-2021"""
$ws.Range("A76").Copy()
$ws.Range("A76").PasteSpecial(-4163)

$row76 = [ordered]@{
    "B"  = 372763
    "C"  = 51306
    "D"  = 29778
    "E"  = 21528
    "F"  = 9855
    "G"  = 2092
    "H"  = 7763
    "I"  = 1130
    "J"  = 11928
    "K"  = 0
    "L"  = 11928
    "M"  = 187813
    "N"  = 167005
    "O"  = 9064
    "P"  = 11745
    "Q"  = 679
    "R"  = 679
    "S"  = 110053
    "T"  = -277794
    "U"  = 650557
    "V"  = 69533
    "W"  = 3
    "X"  = 69530
    "Y"  = 165558
    "Z"  = 5940
    "AA" = 159619
    "AB" = 302534
    "AC" = 102526
    "AD" = 200008
    "AE" = 112932
    "AF" = 6884
    "AG" = 4299
}

foreach ($col in $row76.Keys) {
    $ws.Range("$col`76").Value = $row76[$col]
}
